$d = $word.ActiveDocument

# --- Change 1: merge the title paragraphs -------------------------------
# Before: para 1 = "Uliana Zeigman."   para 2 = "KDD-CUP-98 "
# After:  para 1 = "KDD-CUP-98   final project by  Uliana Zeigman."
#         para 2 = "" (paragraph kept, but its run/text removed)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertBefore("KDD-CUP-98   final project by  ")

$p2 = $d.Paragraphs.Item(2)
$p2r = $p2.Range
$p2NoMark = $d.Range($p2r.Start, $p2r.End - 1)
$p2NoMark.Text = ""

# --- Change 2: insert a new empty paragraph near the end ----------------
# Before end: "...Regression models." / <bold empty para> / <empty para>
# After end:  "...Regression models." / <bold empty para> / <NEW empty para> / <empty para>
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$pLast.Range.InsertParagraphBefore()

# The engine leaves a stray empty run behind on a paragraph split; clean
# it up by forcing a real text mutation (insert a throwaway character,
# then delete it) on the newly created paragraph so it serializes with
# only a <w:pPr> (matching the other blank paragraphs in the document).
$newPara = $d.Paragraphs.Item($count)
$cleanRng = $newPara.Range.Duplicate
$cleanRng.MoveEnd(1, -1)
$cleanRng.InsertBefore("X")
$cleanRng2 = $d.Paragraphs.Item($count).Range.Duplicate
$cleanRng2.MoveEnd(1, -1)
$cleanRng2.Text = ""

Write-Output ("Paragraphs.Count = " + $d.Paragraphs.Count)
